$d = $word.ActiveDocument

# 1. Remove the _GoBack bookmark from its old location (before "Check when syncing...")
$b = $d.Bookmarks.Item("_GoBack")
$b.Delete()

# 2. Insert the new paragraphs (with new check item, spacer paragraph, and new list item)
#    at the very end of the document content, after "Done Dec 15".
$xmlFragment = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
    <w:numPr>
      <w:ilvl w:val="0"/>
      <w:numId w:val="24"/>
    </w:numPr>
  </w:pPr>
  <w:r>
    <w:t xml:space="preserve">Add check </w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve">during sync </w:t>
  </w:r>
  <w:bookmarkStart w:id="0" w:name="_GoBack"/>
  <w:bookmarkEnd w:id="0"/>
  <w:r>
    <w:t xml:space="preserve">to </w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve">make sure all items that DB thinks are in </w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:t>radan</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:t xml:space="preserve"> project are still there.</w:t>
  </w:r>
</w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    <w:ind w:left="360"/>
  </w:pPr>
</w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
    <w:numPr>
      <w:ilvl w:val="0"/>
      <w:numId w:val="24"/>
    </w:numPr>
  </w:pPr>
  <w:r>
    <w:t>Implement drop down menus in ‘add item’ dialog to contain unfinished batches or schedules’</w:t>
  </w:r>
</w:p>
'@

$end = $d.Content.End
$r = $d.Range($end, $end)
$r.InsertXML($xmlFragment)
